# Refresh the cryptocurrency market snapshot on Sheet1 (price + 1h volume change).
# Two coin pairs (37/38, 40/41, 44/45, 50/51) also swapped rank order, so their
# Coin/Link cells are rewritten along with Price/Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "64.720.98"
$ws.Range("E2").Value = "  +3.61%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.463.17"
$ws.Range("E3").Value = "  +4.17%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D5").Value = "578.17"
$ws.Range("E5").Value = "  +4.51%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D6").Value = "156.52"
$ws.Range("E6").Value = "  +3.63%  "

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D7").Value = "0.999"

# Row 8: LidoStakedEther
$ws.Range("D8").Value = "3.465.78"
$ws.Range("E8").Value = "  +4.21%  "

# Row 9: XRP
$ws.Range("D9").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D9").Value = "0.558"
$ws.Range("E9").Value = "  +5.66%  "

# Row 10: Toncoin
$ws.Range("E10").Value = "  +0.86%  "

# Row 11: Dogecoin
$ws.Range("E11").Value = "  +6.74%  "

# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  +2.48%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.064.41"
$ws.Range("E13").Value = "  +4.32%  "

# Row 14: TRON
$ws.Range("E14").Value = "  -1.57%  "

# Row 15: ShibaInu
$ws.Range("E15").Value = "  +10.32%  "

# Row 16: Avalanche
$ws.Range("D16").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D16").Value = "27.94"
$ws.Range("E16").Value = "  +3.71%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "64.770.19"
$ws.Range("E17").Value = "  +3.78%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "3.460.64"
$ws.Range("E18").Value = "  +3.94%  "

# Row 19: Polkadot
$ws.Range("D19").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  -0.92%  "

# Row 20: Chainlink
$ws.Range("D20").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D20").Value = "14.42"
$ws.Range("E20").Value = "  +4.35%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D21").Value = "397.89"
$ws.Range("E21").Value = "  +3.70%  "

# Row 22: Uniswap
$ws.Range("E22").Value = "  +0.84%  "

# Row 23: Polygon
$ws.Range("E23").Value = "  +1.95%  "

# Row 24: Litecoin
$ws.Range("D24").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D24").Value = "73.26"
$ws.Range("E24").Value = "  +3.47%  "

# Row 25: Dai
$ws.Range("E25").Value = "  -0.23%  "

# Row 26: PEPE
$ws.Range("D26").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +26.60%  "

# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D27").Value = "9.62"
$ws.Range("E27").Value = "  +8.88%  "

# Row 28: Kaspa
$ws.Range("E28").Value = "  +2.35%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("D29").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.37%  "

# Row 30: NEARProtocol
$ws.Range("D30").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D30").Value = "6.06"
$ws.Range("E30").Value = "  +9.30%  "

# Row 31: RenderToken
$ws.Range("D31").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D31").Value = "6.79"
$ws.Range("E31").Value = "  +6.50%  "

# Row 32: PancakeSwap
$ws.Range("E32").Value = "  +3.61%  "

# Row 33: Fetch.AI
$ws.Range("E33").Value = "  +5.71%  "

# Row 34: EthereumClassic
$ws.Range("D34").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D34").Value = "23.88"
$ws.Range("E34").Value = "  +3.97%  "

# Row 35: USDe
$ws.Range("E35").Value = "  +0.06%  "

# Row 36: Aptos
$ws.Range("D36").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D36").Value = "7.05"
$ws.Range("E36").Value = "  +4.75%  "

# Row 37: ImmutableX
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D37").Value = "1.49"
$ws.Range("E37").Value = "  -0.39%  "

# Row 38: Monero
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D38").Value = "160.99"
$ws.Range("E38").Value = "  +0.91%  "

# Row 39: Hedera
$ws.Range("D39").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D39").Value = "0.0785"
$ws.Range("E39").Value = "  +7.22%  "

# Row 40: Stacks
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D40").Value = "1.88"
$ws.Range("E40").Value = "  +0.50%  "

# Row 41: EnergySwap
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D41").Value = "27.72"
$ws.Range("E41").Value = "  +2.89%  "

# Row 42: Maker
$ws.Range("D42").Value = "2.910.77"
$ws.Range("E42").Value = "  +2.33%  "

# Row 43: VeChain
$ws.Range("D43").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D43").Value = "0.0325"
$ws.Range("E43").Value = "  +3.24%  "

# Row 44: Filecoin
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D44").Value = "4.44"
$ws.Range("E44").Value = "  +2.46%  "

# Row 45: Mantle
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D45").Value = "0.775"
$ws.Range("E45").Value = "  +3.25%  "

# Row 46: OKB
$ws.Range("D46").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D46").Value = "42.15"
$ws.Range("E46").Value = "  +3.94%  "

# Row 47: InjectiveProtocol
$ws.Range("D47").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D47").Value = "23.81"
$ws.Range("E47").Value = "  +8.24%  "

# Row 48: ONDO
$ws.Range("E48").Value = "  +6.04%  "

# Row 49: dogwifhat
$ws.Range("D49").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D49").Value = "2.20"
$ws.Range("E49").Value = "  +24.23%  "

# Row 50: Cosmos
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D50").Value = "6.57"
$ws.Range("E50").Value = "  +4.34%  "

# Row 51: SuiNetwork
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"  # keep numeric-looking text as text
$ws.Range("D51").Value = "0.858"
$ws.Range("E51").Value = "  +6.07%  "
